$wb = $excel.ActiveWorkbook

# Rename first sheet "LoginData" -> "0-LoginData"
$wsLogin = $wb.Worksheets.Item("LoginData")
$wsLogin.Name = "0-LoginData"

# BasicIdentityTest: change selection from B3 to A1:B1
$wsBasic = $wb.Worksheets.Item("BasicIdentityTest")
$wsBasic.Range("A1:B1").Select()

# AboutMeDesignerTest: no longer the selected tab (tabSelected removed) -- handled implicitly
# because selecting another sheet/range elsewhere moves the active tab.

# Add new worksheet at the end named "7-MyWorkUserTest"
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLast)
$wsNew.Name = "7-MyWorkUserTest"

$wsNew.Range("A1").Value = "folderName"
$wsNew.Range("B1").Value = "اخطار النقل.png"

$wsNew.Range("C8").Select()

# Finally make sure 0-LoginData tab is the selected/active sheet as before (tabSelected=1 on sheet1)
$wsLogin.Activate()
$wsLogin.Range("B6").Select()
